$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B129:I138").Copy()
$ws.Range("Z129").PasteSpecial(-4122)
Write-Host "done"
